# Adding MPA test automation upload file
# Fills in the "Asset Transaction Type" (BWASL, column M) sample values on
# the "Data" sheet for the mass-retirement upload template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# A scratch cell, well outside the template's used range, used as a
# stepping stone so the "250"/"260" codes land in the workbook as genuine
# text (shared-string) cells -- matching how these transaction-type codes
# are stored elsewhere in this template -- rather than being auto-detected
# as numbers.
$scratch = $ws.Range("AA1")

function Set-TextValue($cell, $text) {
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue $ws.Range("M6") "250"
Set-TextValue $ws.Range("M7") "260"
Set-TextValue $ws.Range("M8") "250"
Set-TextValue $ws.Range("M9") "260"
Set-TextValue $ws.Range("M10") "260"
Set-TextValue $ws.Range("M11") "250"
Set-TextValue $ws.Range("M12") "260"
Set-TextValue $ws.Range("M13") "250"
Set-TextValue $ws.Range("M14") "260"
Set-TextValue $ws.Range("M15") "260"

# Row 16 holds a genuine number (not text) for the same field.
$ws.Range("M16").Value = 250
